$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-formatted cells (Price column) to keep their exact original
# string representation (e.g. trailing zeros, double-dot thousand markers)
# instead of Excel auto-converting them to numeric values on assignment.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.797.06'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.295.24'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.98'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.45'
$ws.Range("E6").Value = '  +2.67%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.501'
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  +3.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.05'
$ws.Range("E10").Value = '  +6.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0790'
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("E12").Value = '  +1.81%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.33'
$ws.Range("E13").Value = '  +10.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.87'
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.651.28'
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.244.04'
$ws.Range("E16").Value = '  -2.93%  '
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.829.39'
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.31'
$ws.Range("E19").Value = '  +4.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("E20").Value = '  +2.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0900'
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.91'
$ws.Range("E22").Value = '  +1.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.31'
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.25'
$ws.Range("E24").Value = '  +13.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("E26").Value = '  -0.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.74'
$ws.Range("E27").Value = '  +2.71%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '168.79'
$ws.Range("E28").Value = '  +1.10%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.47'
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.06'
$ws.Range("E30").Value = '  -5.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.16'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.998'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.01'
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.51'
$ws.Range("E34").Value = '  +1.96%  '
$ws.Range("E35").Value = '  -3.11%  '
$ws.Range("E36").Value = '  +2.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0691'
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.83'
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.102'
$ws.Range("E39").Value = '  +1.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.78'
$ws.Range("E40").Value = '  +1.29%  '
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.984.85'
$ws.Range("E42").Value = '  +0.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0289'
$ws.Range("E43").Value = '  +3.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.23'
$ws.Range("E44").Value = '  -2.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.11'
$ws.Range("E45").Value = '  +3.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.55'
$ws.Range("E46").Value = '  -0.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.89'
$ws.Range("E47").Value = '  +1.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.75'
$ws.Range("E48").Value = '  +5.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.518.22'
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("E50").Value = '  +2.16%  '
$ws.Range("E51").Value = '  -1.48%  '
